$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AQ1").Value = "Week 42"
$ws.Range("AQ2").Value = 3.5
$ws.Range("AQ4").Value = 10
$ws.Range("AQ6").Value = 5
$ws.Range("AQ7").Value = 3.5
$ws.Range("AQ9").Value = 10

# Restore the active-cell selection to match the author's saved view
$ws.Range("AN13").Select() | Out-Null
